# "update table from cruise to month"
#
# Insert a new first sheet "biomass_rda_statistics" summarising the
# overall RDA model R-squared / adjusted R-squared for the full model
# and the backward-selected model. The pre-existing "biomass_rda_axis"
# and "biomass_rda_margin" sheets are left untouched and simply shift
# right by one tab position.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() with no args inserts a new sheet immediately before
# the currently active sheet, i.e. at the very front of this workbook.
$statsSheet = $wb.Worksheets.Add()
$statsSheet.Name = "biomass_rda_statistics"

# Header row
$statsSheet.Range("A1").Value = "Model"
$statsSheet.Range("B1").Value = "R.squared"
$statsSheet.Range("C1").Value = "Adj.R.squared"

$headerRange = $statsSheet.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108

# Data rows
$statsSheet.Range("A2").Value = "Full model"
$statsSheet.Range("B2").Value = 0.3730526866937777
$statsSheet.Range("C2").Value = 0.1975074389680355

$statsSheet.Range("A3").Value = "Backward selected"
$statsSheet.Range("B3").Value = 0.3434621809586746
$statsSheet.Range("C3").Value = 0.1919534534875995
